$p = $ppt.ActivePresentation

$oldValue = ": 6369277534"
$newValue = ": asunm10942510"

# Locate the shape/text-frame that holds the "USER ID" value and update just
# the value portion (the ": <value>" run), leaving the "USER ID" label run,
# other paragraphs, and all other runs/formatting untouched.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tf = $sh.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                $full = $tr.Text
                $idx = $full.IndexOf($oldValue)
                if ($idx -ge 0) {
                    $target = $tr.Characters($idx + 1, $oldValue.Length)
                    $target.Text = $newValue
                }
            }
        }
    }
}
